# 🔄 MAJ automatique BRVM via GitHub Actions
# This script rewrites the BRVM recommendations table (sheet "Recommandations")
# rows 2-45 and the Top_YTD sector-progression table (sheet "Top_YTD") rows 2-11
# with the latest automated data pull.

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# Sheet1 (Recommandations) data rows 2-45
$data1 = @(
    @('NEI-CEDA CI', 0, 4, 3715, 935, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES PUBLICS', 0, 8, 3375.12, 110.73, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AUTRES SECTEURS', 0, 4, 2358.84, 593.32, '🟡 Observer', '➖ Neutre'),
    @('BRVM - DISTRIBUTION', 0, 4, 2011.06, 502.15, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TRANSPORT', 0, 4, 1430.7, 354.92, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AGRICULTURE', 0, 4, 1353.35, 337.72, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 677.52, 170.66, '🟡 Observer', '➖ Neutre'),
    @('BRVM - FINANCES', 0, 4, 589.49, 147.33, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 4, 579.36, 144.79, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 4, 570.37, 142.04, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE                 (**)', 0, 2, 533.82, 266.86, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 4, 502.71, 126.2, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 4, 454.18, 112.22, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE          (**)', 0, 2, 444.46, 222.9, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL                    (**)', 0, 2, 440.29, 220.55, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 4, 381.23, 94.97, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE', 0, 1, 264.92, 264.92, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE         (**)', 0, 1, 222.7, 222.7, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL                   (**)', 0, 1, 220.21, 220.21, '🟡 Observer', '➖ Neutre'),
    @('ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)', 4, 0, 29.43, 7.23, '🟢 Achat', '✅ Renforcer'),
    @('SERVAIR ABIDJAN CI (ABJC)', 3, 0, 18.63, 3.92, '🟢 Achat', '✅ Renforcer'),
    @('SOLIBRA CI (SLBC)', 2, 0, 9.69, 7.49, '🟡 Observer', '➖ Neutre'),
    @('SUCRIVOIRE (SCRC)', 2, 1, 5.51, 1.86, '🟡 Observer', '👀 À surveiller'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 1, 0, 4.55, 4.55, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA NG (BOAN)', 1, 0, 3.56, 3.56, '🟡 Observer', '➖ Neutre'),
    @('BERNABE CI (BNBC)', 1, 1, 3.24, -4.11, '🟡 Observer', '👀 À surveiller'),
    @('BANK OF AFRICA BN (BOAB)', 1, 0, 2.76, 2.76, '🟡 Observer', '➖ Neutre'),
    @('CORIS BANK INTERNATIONAL (CBIBF)', 1, 0, 1.6, 1.6, '🟡 Observer', '➖ Neutre'),
    @('SICABLE CI (CABC)', 1, 1, 0.66, -3.89, '🟡 Observer', '👀 À surveiller'),
    @('FILTISAC CI (FTSC)', 1, 1, 0.08, 2.91, '🟡 Observer', '👀 À surveiller'),
    @('TOTAL', 0, 4, 0, 0, '🟡 Observer', '➖ Neutre'),
    @('SMB CI (SMBC)', 0, 1, -1.05, -1.05, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA SENEGAL (BOAS)', 0, 1, -1.76, -1.76, '🟡 Observer', '➖ Neutre'),
    @('ORANGE COTE D''IVOIRE (ORAC)', 0, 1, -1.95, -1.95, '🟡 Observer', '➖ Neutre'),
    @('ONATEL BF (ONTBF)', 1, 1, -2.14, 1.45, '🟡 Observer', '👀 À surveiller'),
    @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 0, 1, -2.42, -2.42, '🟡 Observer', '➖ Neutre'),
    @('ORAGROUP TOGO (ORGT)', 0, 1, -2.72, -2.72, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK COTE D''''IVOIRE (ECOC)', 0, 1, -3.73, -3.73, '🟡 Observer', '➖ Neutre'),
    @('VIVO ENERGY CI (SHEC)', 0, 1, -3.85, -3.85, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI (STAC)', 0, 2, -4.6, -1.87, '🟡 Observer', '➖ Neutre'),
    @('SICOR CI (SICC)', 0, 2, -4.87, -2.37, '🟡 Observer', '➖ Neutre'),
    @('CFAO MOTORS CI (CFAC)', 0, 1, -5.94, -5.94, '🟡 Observer', '➖ Neutre'),
    @('SAFCA CI (SAFC)', 0, 1, -6.09, -6.09, '🟡 Observer', '➖ Neutre'),
    @('UNILEVER CI (UNLC)', 0, 1, -6.97, -6.97, '🟡 Observer', '➖ Neutre')
)

# Sheet2 (Top_YTD) data rows 2-11
$data2 = @(
    @('BRVM - SERVICES PUBLICS', 9615218.42),
    @('NEI-CEDA CI', 1119330.12),
    @('BRVM - AUTRES SECTEURS', 226167.92),
    @('BRVM - DISTRIBUTION', 131905.25),
    @('BRVM - TRANSPORT', 43772.92),
    @('BRVM - AGRICULTURE', 36817.59),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 5165.54),
    @('BRVM - FINANCES', 3644.6),
    @('BRVM - SERVICES FINANCIERS', 3493.6),
    @('BRVM-PRESTIGE', 3363.42)
)

# Write the "Recommandations" sheet (sheet1): columns A..G, rows 2..45
for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $i + 2
    $row = $data1[$i]
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
}

# Write the "Top_YTD" sheet (sheet2): columns A..B, rows 2..11
for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $i + 2
    $row = $data2[$i]
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
}
